$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two additional managers for testing (Donkey and Monkey), entered
# column by column to mirror how the values were typed into the sheet.
$ws.Range("A4").Value = "Donkey"
$ws.Range("A5").Value = "Monkey"

$ws.Range("B4").Value = "T1111111B"
$ws.Range("B5").Value = "T2222222B"

$ws.Range("C4").Value = 35
$ws.Range("C5").Value = 35

$ws.Range("D4").Value = "Married"
$ws.Range("D5").Value = "Married"

$ws.Range("E4").Value = "password"
$ws.Range("E5").Value = "password"

$ws.Range("B4").Select()
